# "added ink! code snippets"
#
# Three cell-content edits on Sheet1:
#   C9  "eosio"            -> "eosio C++"
#   D10 "Solidity, Vyper"  -> "Solidity, Vyper, Fe"
#   B12 (empty)            -> "Self::env().caller"   (ink! equivalent of msg.sender)
#
# Plus: H12's cell style is normalised to match the rest of its row
# (same visual style as e.g. D12 - no fill), and the sheet view is
# scrolled/selected so B13 is the active cell with row 5 at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits -------------------------------------------------------

$ws.Range("C9").Value = "eosio C++"
$ws.Range("D10").Value = "Solidity, Vyper, Fe"
$ws.Range("B12").Value = "Self::env().caller"

# --- Cosmetic style tidy-up on H12 (drop the now-unused "applyFill" xf) --

$ws.Range("H12").Interior.Pattern = -4142

# --- Sheet view: scroll so row 5 is at the top, select B13 --------------

$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("B13").Select()
